$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A174").Value = 173
$ws.Range("B174").Value = 1
$ws.Range("C174").Value = "2024-06-18 15:13:12"
$ws.Range("D174").Value = 200
$ws.Range("E174").Value = 21

$ws.Range("A175").Value = 174
$ws.Range("B175").Value = 2
$ws.Range("C175").Value = "2024-06-18 15:13:13"
$ws.Range("D175").Value = 200
$ws.Range("E175").Value = 1
